$d = $word.ActiveDocument

# --- Paragraph 3: 'Tiruveedhula' / "Rooda's" proofErr splits ---
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="001D2FC2" w:rsidRDefault="00B63F21" w:rsidP="001D2FC2"><w:r><w:t xml:space="preserve">A) </w:t></w:r><w:r w:rsidR="001D2FC2"><w:t xml:space="preserve">Our emulator is based around a DAC card with signal generating capabilities.  The card was recommended to us by Pavan </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D2FC2" w:rsidRPr="001D2FC2"><w:t>Tiruveedhula</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D2FC2"><w:t xml:space="preserve"> from Austin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rooda''s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> lab – they use it for similar purposes.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r3.InsertXML($xml3)

# --- Paragraph 5: 'Giblett' proofErr split (preserve hyperlink) ---
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00F23EFA" w:rsidRDefault="00F23EFA" w:rsidP="001D2FC2"><w:r><w:t>The DAC card was supplied by Strategic Test Corp (</w:t></w:r><w:r w:rsidR="00F34374"><w:t xml:space="preserve">strategic-test.com; </w:t></w:r><w:r><w:t xml:space="preserve">US address: One Boston Place, 26th floor, Boston MA, 20108). Our contact there is Bob </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Giblett</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:hyperlink r:id="rId4" w:history="1"><w:r w:rsidRPr="0029073F"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>bob.giblett@strategic-test.com</w:t></w:r></w:hyperlink><w:r><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r5.InsertXML($xml5)
# InsertXML drops the rStyle on the hyperlink run in this host; restore it.
$hl = $d.Content
$hlFound = $hl.Find.Execute("bob.giblett@strategic-test.com")
if ($hlFound) { $hl.Style = "Hyperlink" }

# --- Paragraph 7: insert calibration-certificate sentence + _GoBack bookmark ---
$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$xml7 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00452F0E" w:rsidRDefault="00F23EFA" w:rsidP="00766057"><w:r><w:t xml:space="preserve">What we purchased was a model UF2e-6022.  It cost us just over $5000 in 2017.  </w:t></w:r><w:r w:rsidR="00673BF3"><w:t xml:space="preserve">Drivers were included. </w:t></w:r><w:r><w:t xml:space="preserve">The quote </w:t></w:r><w:r w:rsidR="000E0DE8"><w:t xml:space="preserve">we used for the purchase is </w:t></w:r><w:r><w:t>in this directory</w:t></w:r><w:r w:rsidR="000E0DE8"><w:t xml:space="preserve"> (UF2e-6022Quote.pdf)</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="00274965"><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t xml:space="preserve">The calibration certificate for our card is also in this directory (UF2e-6022CalibrationCert.pdf). </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>We think that either the appropriate cables come included with the board, or that we had these rolled into the quote we received. You will want to make sure you get cables if you build one of these systems.</w:t></w:r><w:r w:rsidR="00766057"><w:t xml:space="preserve">  It might be that the cable part number is Cab-3f-9m-200.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r7.InsertXML($xml7)

# --- Paragraph 17: 'So' gramStart/gramEnd + drop old _GoBack bookmark location ---
$p17 = $d.Paragraphs.Item(17)
$r17 = $p17.Range
$xml17 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="007629AA" w:rsidRDefault="007629AA" w:rsidP="001D2FC2"><w:r><w:t>D) Once you have the card, you need to install some drivers.  See EmulatorHardwardSetup.pptx for step-by-step instructions.</w:t></w:r><w:r w:rsidR="00696C01"><w:t xml:space="preserve">  SAS Computing has an image of the computer after we installed the OS, Matlab, and the drivers in 2018.</w:t></w:r><w:r w:rsidR="000A2FA8"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00696C01"><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00696C01"><w:t xml:space="preserve"> if we need to do it again, that''s the place to start.</w:t></w:r><w:r w:rsidR="009E0CC3"><w:t xml:space="preserve"> David has the DVD with the drivers f</w:t></w:r><w:r w:rsidR="00691B19"><w:t>or</w:t></w:r><w:r w:rsidR="009E0CC3"><w:t xml:space="preserve"> the UF2e-6022.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r17.InsertXML($xml17)

Write-Output "edit applied"
